$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 234 to make room for the new data point,
# shifting the existing historical rows down by one.
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new day's data.
# Force column A to be treated as plain text so the date string isn't
# auto-converted into a date serial number (matches the other date cells).
$ws.Range("A234").NumberFormat = "@"
$ws.Range("A234").Value = "12.01.2021"
$ws.Range("A234").ClearFormats()
$ws.Range("B234").Value = 2889
$ws.Range("C234").Value = 11608
$ws.Range("D234").Value = 146
$ws.Range("E234").Value = 2080
$ws.Range("F234").Value = 663
$ws.Range("G234").Value = 0
